# Add data for 2021-10-02
# Rename sheet to reflect new "through" date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Through 2021-09-24"

# Row 9 (July)
$ws.Range("U9").Value = 141
$ws.Range("V9").Value = 0.06619999999999999

# Row 10 (August)
$ws.Range("T10").Value = 6
$ws.Range("U10").Value = 153
$ws.Range("V10").Value = 0.0377

# Row 11 (September) - label update
$ws.Range("A11").Value = "September (through 09-24)"
$ws.Range("F11").Value = 33
$ws.Range("G11").Value = 0.0571
$ws.Range("I11").Value = 56
$ws.Range("J11").Value = 0.0667
$ws.Range("L11").Value = 43
$ws.Range("M11").Value = 0.0851
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 54
$ws.Range("P11").Value = 0.1
$ws.Range("R11").Value = 90
$ws.Range("S11").Value = 0.0323
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = 141
$ws.Range("V11").Value = 0.007

# Row 12 (Total)
$ws.Range("F12").Value = 373
$ws.Range("G12").Value = 0.1034
$ws.Range("I12").Value = 562
$ws.Range("J12").Value = 0.08019999999999999
$ws.Range("L12").Value = 476
$ws.Range("M12").Value = 0.1136
$ws.Range("N12").Value = 42
$ws.Range("O12").Value = 367
$ws.Range("P12").Value = 0.1027
$ws.Range("R12").Value = 826
$ws.Range("S12").Value = 0.0582
$ws.Range("T12").Value = 75
$ws.Range("U12").Value = 1137
$ws.Range("V12").Value = 0.0619
